$d = $word.ActiveDocument

# The paragraph we are targeting is the second "DIFFERENCES" block (the one that
# begins "DIFFERENCES - The information in the Official Ministry record is
# different...". We search on the long, unique sentence that follows "The" so we
# don't collide with the other "DIFFERENCES" occurrence earlier in the document.

$oldText = "The information in the Official Ministry record is different than the record you have submitted. PLEASE FAX the PEN Coordinator any legal documentation that indicates the information in the record you have submitted is correct. Also, listed in this section are any 'changed' PENs.  Please replace your current PEN section with the Official Ministry PEN."

$newText = "The existing student demographic information in PEN is different than the student demographic record^p" + `
  "that you have submitted. Please provide the PEN Coordinator with the current valid Legal Identification indicating^p" + `
  "that the information in the record you have submitted is correct. If the report identifies a change in PEN, Please^p" + `
  "replace your current PEN with the Official Ministry PEN in your school database. MyEducation users NOTE: You must^p" + `
  "have the old PEN nulled by your MyEd HelpDesk."

$range = $d.Content
$found = $range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

if (-not $found) {
    throw "Could not find the target paragraph text to replace."
}

# The replacement above spans several new paragraphs created via the "^p"
# markers. The very last run created this way (the one holding "HelpDesk.")
# ends up at the very end of the document's content with no following
# sibling run to inherit run-formatting from, so it loses its explicit
# Courier New run formatting. Re-apply the same character formatting used
# throughout this section to keep it consistent.

$lastRange = $d.Content
$lastFound = $lastRange.Find.Execute("have the old PEN nulled by your MyEd HelpDesk.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($lastFound) {
    $lastRange.Font.NameAscii = "Courier New"
    $lastRange.Font.NameFarEast = "Times New Roman"
    $lastRange.Font.NameOther = "Courier New"
    $lastRange.Font.Size = 10.5
    $lastRange.Font.Kerning = 0
}

Write-Output "Replaced: $found"
